$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the date values (serial 42078 = 2015-03-15) for B6, C6, E6.
# D6 is left blank, matching the target diff.
$ws.Range("B6").Value2 = 42078
$ws.Range("C6").Value2 = 42078
$ws.Range("E6").Value2 = 42078

# Update the selection to D6, as in the diff.
$ws.Range("D6").Select()
